# Apply the three changes described by the diff:
#  1) Footer date text "30.07.2025" -> "31.07.2025" on slide layout 2
#     (shape "Rectangle 7" holding a datetime field).
#  2) Same footer date change on slide layout 5 ("Rectangle 7").
#  3) On slide 2, inside the grouped shape, the "Textfeld 4" textbox:
#       - shrink its height (cy 577081 -> 415498 EMU)
#       - change its second paragraph's text from
#         "Pauses and resumes the game" to "Pause / Resume"

$p = $ppt.ActivePresentation

# --- 1 & 2: update the date footer text boxes on the slide layouts ---
$master = $p.SlideMaster

$layout2 = $master.CustomLayouts.Item(2)
$dateShape2 = $layout2.Shapes.Item(7)
$dateShape2.TextFrame.TextRange.Text = "31.07.2025"

$layout5 = $master.CustomLayouts.Item(5)
$dateShape5 = $layout5.Shapes.Item(9)
$dateShape5.TextFrame.TextRange.Text = "31.07.2025"

# --- 3: update the "Pause / Resume" button caption on slide 2 ---
$slide2 = $p.Slides.Item(2)
$buttonGroup = $slide2.Shapes.Item(4)
$textBox = $buttonGroup.GroupItems.Item(2)

$secondPara = $textBox.TextFrame.TextRange.Paragraphs(2)
# First retype with text that shares no characters with the old caption so
# the engine replaces the whole run cleanly (keeping the original run
# formatting) instead of diff-splitting into several runs; then set the
# real final caption on top of that.
$secondPara.Text = "ZzZ9999Qq"
$textBox.TextFrame.TextRange.Paragraphs(2).Text = "Pause / Resume"
